$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "failureReason"
$ws.Range("A2").Value = "CTX456161"
$ws.Range("B2").Value = " Rajesh"
$ws.Range("C2").Value = "Shankar"
